$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Withdraw the key from TanLin (was row 7 / key #6). Remove her name and
# close the gap by moving everyone below her up one row in the name
# column (B) only -- the key-number column (A) stays exactly as-is.
$ws.Range("B7").Value  = "丁昊"
$ws.Range("B8").Value  = "李娜"
$ws.Range("B9").Value  = "姜善宸"
$ws.Range("B10").Value = "郑自强"
$ws.Range("B11").Value = "张少永"
$ws.Range("B12").Value = "卢婧宇"
$ws.Range("B13").Value = "杜昂昂"
$ws.Range("B14").ClearContents()

# Update the running totals in the header note: one fewer key in use,
# one more sitting idle (19 total; use 13->12, idle 6->7).
# "（目前共19把，使用13把，闲置6把）"
#   "13" is at characters 11-12, "6" is at character 17.
$ws.Range("E1").Characters(11, 2).Text = "12"
$ws.Range("E1").Characters(17, 1).Text = "7"
